$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 205 (1-based), shifting existing rows 205..297 down to 206..298
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record's data
$ws.Cells.Item(205, 1).Value = 5
$ws.Cells.Item(205, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(205, 3).Value = "Maule"
$ws.Cells.Item(205, 4).Value = 44784
$ws.Cells.Item(205, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(205, 5).Value = 7
$ws.Cells.Item(205, 6).Value = 100112008
$ws.Cells.Item(205, 7).Value = "Coliflor"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 3000
$ws.Cells.Item(205, 11).Value = 1000
$ws.Cells.Item(205, 12).Value = 1000
$ws.Cells.Item(205, 13).Value = 1000
$ws.Cells.Item(205, 14).Value = "`$/unidad"
$ws.Cells.Item(205, 15).Value = "Región del Maule"
$ws.Cells.Item(205, 16).Value = 1000
$ws.Cells.Item(205, 17).Value = 1
$ws.Cells.Item(205, 18).Value = "Hortaliza"
